$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$helper = $ws.Range("ZZ1")

$ws.Range("D2").Value = "30.244.47"
$ws.Range("E2").Value = "  +1.19%  "
$ws.Range("D3").Value = "2.087.05"
$ws.Range("E3").Value = "  -1.39%  "
$helper.Value = "'1.005"
$helper.Copy()
$ws.Range("D4").PasteSpecial(-4163)
$ws.Range("E4").Value = "  -0.23%  "
$helper.Value = "'341.12"
$helper.Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = "  -2.01%  "
$helper.Value = "'1.004"
$helper.Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = "  -0.16%  "
$helper.Value = "'0.5282"
$helper.Copy()
$ws.Range("D7").PasteSpecial(-4163)
$ws.Range("E7").Value = "  +1.96%  "
$helper.Value = "'0.4385"
$helper.Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("E8").Value = "  -1.80%  "
$helper.Value = "'54.95"
$helper.Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("E9").Value = "  +1.94%  "
$helper.Value = "'0.09343"
$helper.Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Value = "  -0.09%  "
$helper.Value = "'1.174"
$helper.Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Value = "  -0.77%  "
$helper.Value = "'24.53"
$helper.Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("E12").Value = "  -2.81%  "
$helper.Value = "'8.495"
$helper.Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("E13").Value = "  +2.16%  "
$helper.Value = "'6.862"
$helper.Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Value = "  +0.30%  "
$ws.Range("D15").Value = "2.073.15"
$ws.Range("E15").Value = "  -1.21%  "
$helper.Value = "'101.49"
$helper.Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("E16").Value = "  -1.15%  "
$helper.Value = "'0.00001157"
$helper.Copy()
$ws.Range("D17").PasteSpecial(-4163)
$ws.Range("E17").Value = "  -0.68%  "
$helper.Value = "'1.005"
$helper.Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("E18").Value = "  -0.21%  "
$helper.Value = "'21.02"
$helper.Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Value = "  -2.29%  "
$helper.Value = "'0.06710"
$helper.Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Value = "  +0.58%  "
$helper.Value = "'6.280"
$helper.Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = "  -0.46%  "
$helper.Value = "'1.003"
$helper.Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E22").Value = "  -0.25%  "
$ws.Range("D23").Value = "30.262.05"
$ws.Range("E23").Value = "  +1.12%  "
$helper.Value = "'12.38"
$helper.Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = "  -2.81%  "
$helper.Value = "'2.326"
$helper.Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("E26").Value = "  -1.87%  "
$helper.Value = "'6.841"
$helper.Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = "  +7.00%  "
$helper.Value = "'162.31"
$helper.Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E28").Value = "  -0.11%  "
$helper.Value = "'2.481"
$helper.Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Value = "  -2.93%  "
$helper.Value = "'133.43"
$helper.Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Value = "  -0.44%  "
$helper.Value = "'1.127"
$helper.Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Value = "  -2.13%  "
$helper.Value = "'0.1048"
$helper.Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Value = "  -0.83%  "
$helper.Value = "'1.660"
$helper.Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("E33").Value = "  -7.34%  "
$helper.Value = "'6.253"
$helper.Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("E34").Value = "  +0.00%  "
$helper.Value = "'3.903"
$helper.Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value = "  -1.61%  "
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$helper.Value = "'0.02607"
$helper.Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Value = "  +0.47%  "
$ws.Range("B37").Value = "FraxShare"
$ws.Range("C37").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$helper.Value = "'9.916"
$helper.Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Value = "  -8.52%  "
$helper.Value = "'0.06733"
$helper.Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Value = "  -1.10%  "
$ws.Range("B39").Value = "Aptos"
$ws.Range("C39").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$helper.Value = "'12.56"
$helper.Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Value = "  -1.37%  "
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$helper.Value = "'1.342"
$helper.Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Value = "  -0.10%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$helper.Value = "'0.6934"
$helper.Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value = "  -1.31%  "
$helper.Value = "'0.2198"
$helper.Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Value = "  -2.02%  "
$helper.Value = "'0.6732"
$helper.Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Value = "  -1.90%  "
$helper.Value = "'2.380"
$helper.Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Value = "  +0.72%  "
$helper.Value = "'14.24"
$helper.Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Value = "  -2.38%  "
$helper.Value = "'1.004"
$helper.Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = "  -0.14%  "
$helper.Value = "'1.289"
$helper.Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Value = "  +5.73%  "
$ws.Range("E48").Value = "  -0.09%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$helper.Value = "'0.00000000344"
$helper.Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("E49").Value = "  -2.48%  "
$ws.Range("B50").Value = "ThetaToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$helper.Value = "'1.207"
$helper.Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("E50").Value = "  +1.75%  "
$helper.Value = "'1.208"
$helper.Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value = "  -1.47%  "
$helper.Clear()
$excel.CutCopyMode = $false
